# Apply crypto price/volume updates per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.935.16'
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('E3').Value = '  -0.50%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.86%  '
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('E10').Value = '  -1.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0795'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.35%  '
$ws.Range('D12').Value = '1.864.48'
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.26'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.17%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.605.46'
$ws.Range('E14').Value = '  -2.02%  '
$ws.Range('E15').Value = '  -1.13%  '
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.92'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('D18').Value = '25.943.68'
$ws.Range('E18').Value = '  -1.03%  '
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '192.99'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.17%  '
$ws.Range('E21').Value = '  -1.49%  '
$ws.Range('E22').Value = '  -1.21%  '
$ws.Range('E23').Value = '  -0.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '144.10'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.50%  '
$ws.Range('E25').Value = '  +0.81%  '
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.130'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.95%  '
$ws.Range('B27').Value = 'BinanceUSD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.21%  '
$ws.Range('E28').Value = '  -1.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.55'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.25'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.52%  '
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('E32').Value = '  -1.12%  '
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('E34').Value = '  -4.18%  '
$ws.Range('E35').Value = '  +1.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.901'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.28%  '
$ws.Range('D37').Value = '1.138.34'
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('E38').Value = '  -1.64%  '
$ws.Range('E39').Value = '  -1.22%  '
$ws.Range('E40').Value = '  +0.30%  '
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.49'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.25'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.92%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.798'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('D45').Value = '1.774.31'
$ws.Range('E45').Value = '  -0.46%  '
$ws.Range('E46').Value = '  +2.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.62'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.56%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0532'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.47'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.68'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.39%  '
$ws.Range('E51').Value = '  -0.83%  '
